# Add two new columns, I (I0) and J (IF), to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---------------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, centered, bordered) from H1 onto the two
# new header cells so they match the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data rows -----------------------------------------------------------
# For rows 2-19, I is always 1 and J duplicates the value already present
# in column H for that row.
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 8).Value2
}

# Row 20 is a special case that does not follow the pattern above.
$ws.Range("I20").Value = 5
$ws.Range("J20").Value = 6
